$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. Fix the typo in row 5's time range: 14.00-16.00 -> 14.00-15.00
$ws.Range("B5").Value = "9.00-10.45, 11:45-13.00, 14.00-15.00"

# Hours logged on 29 syys corrected from 5 to 4
$ws.Range("G5").Value = 4

# 2. Insert a new diary entry for "1 loka" in row 6
$ws.Range("A6").Value = "1 loka"
$ws.Range("B6").Value = "17.50-18.50, 19.15-20.45"
$ws.Range("C6").Value = "Edelleen partikkelia yritän saada näkymään. Projektirakenteen siistimistä, C++ syntaksin ihmettelyä."
$ws.Range("D6").Value = "Oppiminen olisi tehokkaampaa, jos olisi enemmän sujut kielen kanssa. Kyllä tämä tästä jossakin vaiheessa."
$ws.Range("E6").Value = "Pakko sanoa, että luulin jo päässeeni sinuiksi enemmän tällaisten aloittelijan virheiden kanssa, mutta se on hyvä opettaja kun koodi ei toimi. Täytyy huolehtia, että jää aikaa opetella myös asiasisältöä, mutta tässä toistaiseksi ollut ihan hommaa että saa omat solmut availtua."
$ws.Range("F6").Value = "Eli näköjään kun käyttä glad, ja glfw tai muuta kirjastoa, tulee aina ensin sisällyttää glad, riippumatta tarvitseeko sitä ko tiedostossa. Tällöin vältytään include guard virheeltä joka tulkitsee gl.h kirjaston sisällyttämisen kahdesti."
$ws.Range("G6").Value = 2.5

# Match formatting of the row above (wrap text on content cells, time format on Kello, row height)
$ws.Range("B6").NumberFormat = $ws.Range("B5").NumberFormat
$ws.Range("B6").WrapText = $true
$ws.Range("C6:F6").WrapText = $true
$ws.Rows.Item(6).RowHeight = 116

# 3. Widen META / Tunnit columns to fit the new content
$ws.Columns.Item(6).ColumnWidth = 31.5
$ws.Columns.Item(7).ColumnWidth = 14.3

# 4. Update the active selection to reflect where editing finished
$ws.Range("G7").Select()

$wb.Save()
